$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $cellRef, $text) {
    $range = $sheet.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '30.296.51'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.929.79'
$ws.Range('E3').Value = '  -0.64%  '
Set-TextCell $ws 'D4' '0.9993'
$ws.Range('E4').Value = '  -0.20%  '
Set-TextCell $ws 'D5' '0.7462'
$ws.Range('E5').Value = '  +3.49%  '
Set-TextCell $ws 'D6' '249.71'
$ws.Range('E6').Value = '  -0.85%  '
Set-TextCell $ws 'D7' '0.9998'
$ws.Range('E7').Value = '  -0.11%  '
Set-TextCell $ws 'D8' '0.3226'
$ws.Range('E8').Value = '  -3.46%  '
Set-TextCell $ws 'D9' '27.91'
$ws.Range('E9').Value = '  -3.00%  '
Set-TextCell $ws 'D10' '0.07113'
$ws.Range('E10').Value = '  -4.06%  '
Set-TextCell $ws 'D11' '0.7896'
$ws.Range('E11').Value = '  -3.34%  '
Set-TextCell $ws 'D12' '0.08022'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').Value = '1.934.25'
$ws.Range('E13').Value = '  -0.39%  '
Set-TextCell $ws 'D14' '5.393'
$ws.Range('E14').Value = '  -1.78%  '
Set-TextCell $ws 'D15' '94.63'
$ws.Range('E15').Value = '  -0.64%  '
Set-TextCell $ws 'D16' '14.57'
$ws.Range('E16').Value = '  -2.51%  '
$ws.Range('D17').Value = '30.311.06'
$ws.Range('E17').Value = '  -0.24%  '
Set-TextCell $ws 'D18' '252.68'
$ws.Range('E18').Value = '  -0.48%  '
Set-TextCell $ws 'D19' '0.000008070'
$ws.Range('E19').Value = '  -4.92%  '
Set-TextCell $ws 'D20' '5.750'
$ws.Range('E20').Value = '  -2.47%  '
$ws.Range('D21').Value = '2.185.89'
Set-TextCell $ws 'D22' '0.9997'
$ws.Range('E22').Value = '  -0.09%  '
Set-TextCell $ws 'D23' '0.9992'
$ws.Range('E23').Value = '  -0.10%  '
Set-TextCell $ws 'D24' '6.838'
$ws.Range('E24').Value = '  -2.14%  '
Set-TextCell $ws 'D25' '9.584'
$ws.Range('E25').Value = '  -2.89%  '
Set-TextCell $ws 'D26' '164.24'
$ws.Range('E26').Value = '  +0.72%  '
Set-TextCell $ws 'D27' '19.09'
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D28' '0.1336'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws 'D29' '2.277'
$ws.Range('E29').Value = '  -5.75%  '
Set-TextCell $ws 'D30' '1.360'
$ws.Range('E30').Value = '  +1.06%  '
Set-TextCell $ws 'D31' '1.535'
$ws.Range('E31').Value = '  -2.50%  '
Set-TextCell $ws 'D32' '4.413'
$ws.Range('E32').Value = '  -1.14%  '
Set-TextCell $ws 'D33' '4.156'
$ws.Range('E33').Value = '  -2.44%  '
Set-TextCell $ws 'D34' '0.05115'
$ws.Range('E34').Value = '  -3.50%  '
Set-TextCell $ws 'D35' '1.294'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('E36').Value = '  -1.45%  '
Set-TextCell $ws 'D38' '0.01976'
$ws.Range('E38').Value = '  -0.83%  '
Set-TextCell $ws 'D39' '2.799'
$ws.Range('E39').Value = '  -1.82%  '
Set-TextCell $ws 'D40' '78.07'
$ws.Range('E40').Value = '  -3.96%  '
Set-TextCell $ws 'D41' '6.405'
$ws.Range('E41').Value = '  -3.25%  '
Set-TextCell $ws 'D42' '0.4512'
$ws.Range('E42').Value = '  -1.36%  '
Set-TextCell $ws 'D43' '1.989'
$ws.Range('E43').Value = '  -3.03%  '
Set-TextCell $ws 'D44' '0.8457'
$ws.Range('E44').Value = '  -0.31%  '
Set-TextCell $ws 'D45' '0.9998'
$ws.Range('E45').Value = '  -0.11%  '
Set-TextCell $ws 'D46' '101.48'
$ws.Range('E46').Value = '  -1.44%  '
Set-TextCell $ws 'D47' '7.540'
$ws.Range('E47').Value = '  +0.34%  '
Set-TextCell $ws 'D48' '9.790'
$ws.Range('E48').Value = '  -0.88%  '
Set-TextCell $ws 'D49' '984.23'
$ws.Range('E49').Value = '  +11.86%  '
Set-TextCell $ws 'D50' '37.09'
$ws.Range('E50').Value = '  +0.45%  '
Set-TextCell $ws 'D51' '0.06076'
